$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1531.4375
$ws.Range("I39").Value = 576.8182
$ws.Range("K39").Value = 1730.4546
$ws.Range("M39").Value = -1434.4546

$ws.Range("H64").Value = 5498.5
$ws.Range("I64").Value = 5495
$ws.Range("J64").Value = 5499.375
$ws.Range("K64").Value = 5495
$ws.Range("L64").Value = 5499.375
$ws.Range("M64").Value = -5247
$ws.Range("N64").Value = -5995.375

$ws.Range("H67").Value = 5498.5
$ws.Range("I67").Value = 5495
$ws.Range("J67").Value = 5499.375
$ws.Range("K67").Value = 5495
$ws.Range("L67").Value = 5499.375
$ws.Range("M67").Value = -4637
$ws.Range("N67").Value = -7215.375

$ws.Range("H100").Value = 3522.7942
$ws.Range("I100").Value = 3396
$ws.Range("J100").Value = 4114.5
$ws.Range("K100").Value = 3396
$ws.Range("L100").Value = 4114.5
$ws.Range("M100").Value = -2855
$ws.Range("N100").Value = -5196.5

$ws.Range("H107").Value = 939.38464
$ws.Range("I107").Value = 705.9091
$ws.Range("J107").Value = 2223.5
$ws.Range("K107").Value = 705.9091
$ws.Range("L107").Value = 2223.5
$ws.Range("M107").Value = 1214.0909
$ws.Range("N107").Value = -6063.5

$ws.Range("H116").Value = 15781.27
$ws.Range("I116").Value = 13960.25
$ws.Range("J116").Value = 16590.611
$ws.Range("K116").Value = 13960.25
$ws.Range("L116").Value = 16590.611
$ws.Range("M116").Value = -10518.25
$ws.Range("N116").Value = -23474.611

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 16070.8
$ws.Range("I61").Value = 18774.875
$ws.Range("J61").Value = 5254.5
$ws.Range("K61").Value = 18774.875
$ws.Range("L61").Value = 5254.5
$ws.Range("M61").Value = -18562.875
$ws.Range("N61").Value = -5678.5

$ws.Range("H110").Value = 2121.353
$ws.Range("I110").Value = 1051.5454
$ws.Range("K110").Value = 1051.5454
$ws.Range("M110").Value = 993.4546

$ws.Range("H123").Value = 78636.37
$ws.Range("I123").Value = 78636.37
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 78636.37
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -73736.37
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 28951.37
$ws.Range("I132").Value = 37928.75
$ws.Range("K132").Value = 113786.25
$ws.Range("M132").Value = -111256.25

$ws.Range("H135").Value = 49998.668
$ws.Range("J135").Value = 49998.668
$ws.Range("L135").Value = 49998.668
$ws.Range("N135").Value = -60138.668

$ws.Range("H136").Value = 16070.8
$ws.Range("I136").Value = 18774.875
$ws.Range("J136").Value = 5254.5
$ws.Range("K136").Value = 56324.625
$ws.Range("L136").Value = 15763.5
$ws.Range("M136").Value = -53774.625
$ws.Range("N136").Value = -20863.5

$ws.Range("H140").Value = 79800
$ws.Range("J140").Value = 79800
$ws.Range("L140").Value = 79800
$ws.Range("N140").Value = -90160

$ws.Range("H141").Value = 149900
$ws.Range("J141").Value = 149900
$ws.Range("L141").Value = 149900
$ws.Range("N141").Value = -160260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7695062
$ws.Range("I20").Value = 12502368
$ws.Range("J20").Value = 3371.6
$ws.Range("K20").Value = 12502368
$ws.Range("L20").Value = 3371.6
$ws.Range("M20").Value = -12502121
$ws.Range("N20").Value = -3865.6

$ws.Range("H33").Value = 3021
$ws.Range("I33").Value = 3021
$ws.Range("K33").Value = 3021
$ws.Range("M33").Value = -2685

$ws.Range("H105").Value = 3285.5945
$ws.Range("I105").Value = 3118.2812
$ws.Range("K105").Value = 3118.2812
$ws.Range("M105").Value = -1371.2812

$ws.Range("H107").Value = 3107.28
$ws.Range("I107").Value = 2043.3125
$ws.Range("J107").Value = 4998.778
$ws.Range("K107").Value = 2043.3125
$ws.Range("L107").Value = 4998.778
$ws.Range("M107").Value = -123.3125
$ws.Range("N107").Value = -8838.778

$ws.Range("H137").Value = 44999.5
$ws.Range("J137").Value = 44999.5
$ws.Range("L137").Value = 44999.5
$ws.Range("N137").Value = -55199.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 247.07143
$ws.Range("J7").Value = 407.25
$ws.Range("L7").Value = 407.25
$ws.Range("N7").Value = -633.25

$ws.Range("H31").Value = 1554.4445
$ws.Range("I31").Value = 1627.7693
$ws.Range("J31").Value = 1363.8
$ws.Range("K31").Value = 1627.7693
$ws.Range("L31").Value = 1363.8
$ws.Range("M31").Value = -1332.7693
$ws.Range("N31").Value = -1953.8

$ws.Range("H34").Value = 1554.4445
$ws.Range("I34").Value = 1627.7693
$ws.Range("J34").Value = 1363.8
$ws.Range("K34").Value = 1627.7693
$ws.Range("L34").Value = 1363.8
$ws.Range("M34").Value = -1425.7693
$ws.Range("N34").Value = -1767.8

$ws.Range("H35").Value = 1071.75
$ws.Range("I35").Value = 1095.6666
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 1095.6666
$ws.Range("L35").Value = 1000
$ws.Range("M35").Value = -801.6666
$ws.Range("N35").Value = -1588

$ws.Range("H105").Value = 1523.1538
$ws.Range("I105").Value = 1555.64
$ws.Range("K105").Value = 1555.64
$ws.Range("M105").Value = 191.3599999999999

$ws.Range("H107").Value = 3024.3333
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 3832.4443
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 3832.4443
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -7672.4443

$ws.Range("H122").Value = 2385.7856
$ws.Range("J122").Value = 2294
$ws.Range("L122").Value = 6882
$ws.Range("N122").Value = -11782

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1032.25
$ws.Range("I109").Value = 1032.25
$ws.Range("K109").Value = 3096.75
$ws.Range("M109").Value = -2056.75

$ws.Range("H125").Value = 7264.5
$ws.Range("I125").Value = 7264.5
$ws.Range("K125").Value = 21793.5
$ws.Range("M125").Value = -16873.5

$ws.Range("H129").Value = 3224.3684
$ws.Range("I129").Value = 2259.8572
$ws.Range("J129").Value = 3787
$ws.Range("K129").Value = 6779.571599999999
$ws.Range("L129").Value = 11361
$ws.Range("M129").Value = -1779.571599999999
$ws.Range("N129").Value = -21361

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3339.6667
$ws.Range("I122").Value = 2278.2307
$ws.Range("K122").Value = 6834.6921
$ws.Range("M122").Value = -4384.6921

$ws.Range("H132").Value = 114638.11
$ws.Range("I132").Value = 128467.875
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 385403.625
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -382873.625
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3928.25
$ws.Range("I32").Value = 3928.25
$ws.Range("K32").Value = 3928.25
$ws.Range("M32").Value = -3611.25

$ws.Range("H40").Value = 7497.5835
$ws.Range("I40").Value = 6775.778
$ws.Range("J40").Value = 9663
$ws.Range("K40").Value = 6775.778
$ws.Range("L40").Value = 9663
$ws.Range("M40").Value = -6639.778
$ws.Range("N40").Value = -9935

$ws.Range("H46").Value = 24673.4
$ws.Range("I46").Value = 29686.125
$ws.Range("J46").Value = 4622.5
$ws.Range("K46").Value = 29686.125
$ws.Range("L46").Value = 4622.5
$ws.Range("M46").Value = -29498.125
$ws.Range("N46").Value = -4998.5

$ws.Range("H55").Value = 643.3333
$ws.Range("I55").Value = 390.9091
$ws.Range("K55").Value = 390.9091
$ws.Range("M55").Value = -217.9091

$ws.Range("H58").Value = 33665
$ws.Range("I58").Value = 33665
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 33665
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -33405
$ws.Range("N58").ClearContents()

$ws.Range("H62").Value = 27165.75
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 27165.75
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 27165.75
$ws.Range("N62").Value = -28413.75
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 27165.75
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 27165.75
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 81497.25
$ws.Range("N65").Value = -87737.25
$ws.Range("M65").ClearContents()

$ws.Range("H68").Value = 6022
$ws.Range("I68").Value = 4141.6665
$ws.Range("J68").Value = 8278.4
$ws.Range("K68").Value = 4141.6665
$ws.Range("L68").Value = 8278.4
$ws.Range("M68").Value = -3392.6665
$ws.Range("N68").Value = -9776.4

$ws.Range("H71").Value = 6022
$ws.Range("I71").Value = 4141.6665
$ws.Range("J71").Value = 8278.4
$ws.Range("K71").Value = 20708.3325
$ws.Range("L71").Value = 41392
$ws.Range("M71").Value = -16964.3325
$ws.Range("N71").Value = -48880

$ws.Range("H100").Value = 1697.5
$ws.Range("I100").Value = 1084
$ws.Range("J100").Value = 3538
$ws.Range("K100").Value = 1084
$ws.Range("L100").Value = 3538
$ws.Range("M100").Value = -543
$ws.Range("N100").Value = -4620

$ws.Range("H132").Value = 65324.05
$ws.Range("I132").Value = 75675.94
$ws.Range("J132").Value = 6663.3335
$ws.Range("K132").Value = 227027.82
$ws.Range("L132").Value = 19990.0005
$ws.Range("M132").Value = -224497.82
$ws.Range("N132").Value = -25050.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1292.75
$ws.Range("I107").Value = 756.2
$ws.Range("K107").Value = 2268.6
$ws.Range("M107").Value = -348.6000000000004

$ws.Range("H121").Value = 59918
$ws.Range("J121").Value = 59918
$ws.Range("L121").Value = 59918
$ws.Range("N121").Value = -63412

$ws.Range("H122").Value = 1488.125
$ws.Range("I122").Value = 1488.125
$ws.Range("K122").Value = 4464.375
$ws.Range("M122").Value = -2014.375

$ws.Range("H126").Value = 36922.324
$ws.Range("I126").Value = 42983.81
$ws.Range("K126").Value = 128951.43
$ws.Range("M126").Value = -126481.43

$ws.Range("H132").Value = 25410.605
$ws.Range("I132").Value = 25944.191
$ws.Range("K132").Value = 77832.573
$ws.Range("M132").Value = -75302.573

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
